$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Gnai2","Adora1","M2",3,1,150.0354306666667,450.106292,0.4152507364956075,0.4152507364956075,2,0.6666666666666666,0.006408333333333333,0.019225,0.01109363955048399,0.01109363955048399,0.9614770515222221,8.653293463699999,0.004606641993755278,0.004606641993755278),
    @("ECs","Gnai2","Adora1","sCs",3,1,150.0354306666667,450.106292,0.4152507364956075,0.4152507364956075,3,1,0.57125,1.71375,0.988906360449516,0.988906360449516,85.70773976833333,771.3696579150001,0.4106440945018522,0.4106440945018521),
    @("FAPs","Gnai2","Adora1","M2",3,1,68.382243,205.146729,0.1892604742946246,0.1892604742946246,2,0.6666666666666666,0.006408333333333333,0.019225,0.01109363955048399,0.01109363955048399,0.438216207225,3.943945865025,0.002099587482978207,0.002099587482978206),
    @("FAPs","Gnai2","Adora1","sCs",3,1,68.382243,205.146729,0.1892604742946246,0.1892604742946246,3,1,0.57125,1.71375,0.988906360449516,0.988906360449516,39.06335631375001,351.57020682375,0.1871608868116464,0.1871608868116464),
    @("M2","Gnai2","Adora1","M2",3,1,104.737245,314.211735,0.2898796499701289,0.2898796499701289,2,0.6666666666666666,0.006408333333333333,0.019225,0.01109363955048399,0.01109363955048399,0.6711911783749999,6.040720605374999,0.003215820349789078,0.003215820349789078),
    @("M2","Gnai2","Adora1","sCs",3,1,104.737245,314.211735,0.2898796499701289,0.2898796499701289,3,1,0.57125,1.71375,0.988906360449516,0.988906360449516,59.83115120625,538.4803608562499,0.2866638296203399,0.2866638296203398),
    @("sCs","Gnai2","Adora1","M2",3,1,38.15794,114.47382,0.105609139239639,0.105609139239639,2,0.6666666666666666,0.006408333333333333,0.019225,0.01109363955048399,0.01109363955048399,0.2445287988333333,2.2007591895,0.00117158972396143,0.00117158972396143),
    @("sCs","Gnai2","Adora1","sCs",3,1,38.15794,114.47382,0.105609139239639,0.105609139239639,3,1,0.57125,1.71375,0.988906360449516,0.988906360449516,21.797723225,196.179509025,0.1044375495156776,0.1044375495156775)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
